$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Alcam"
$ws.Range("C2").Value = "Cd6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 62.12558000000001
$ws.Range("H2").Value = 186.37674
$ws.Range("I2").Value = 0.9736910227596813
$ws.Range("J2").Value = 0.9736910227596813
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1475986666666667
$ws.Range("N2").Value = 0.442796
$ws.Range("O2").Value = 0.6430044319495352
$ws.Range("P2").Value = 0.643004431949535
$ws.Range("Q2").Value = 9.169652773893336
$ws.Range("R2").Value = 82.52687496504001
$ws.Range("S2").Value = 0.6260876429839508
$ws.Range("T2").Value = 0.6260876429839507

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Alcam"
$ws.Range("C3").Value = "Cd6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 62.12558000000001
$ws.Range("H3").Value = 186.37674
$ws.Range("I3").Value = 0.9736910227596813
$ws.Range("J3").Value = 0.9736910227596813
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.08194666666666665
$ws.Range("N3").Value = 0.24584
$ws.Range("O3").Value = 0.356995568050465
$ws.Range("P3").Value = 0.356995568050465
$ws.Range("Q3").Value = 5.090984195733333
$ws.Range("R3").Value = 45.8188577616
$ws.Range("S3").Value = 0.3476033797757306
$ws.Range("T3").Value = 0.3476033797757306

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Alcam"
$ws.Range("C4").Value = "Cd6"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5683613333333334
$ws.Range("H4").Value = 1.705084
$ws.Range("I4").Value = 0.008907897969731461
$ws.Range("J4").Value = 0.008907897969731461
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1475986666666667
$ws.Range("N4").Value = 0.442796
$ws.Range("O4").Value = 0.6430044319495352
$ws.Range("P4").Value = 0.643004431949535
$ws.Range("Q4").Value = 0.08388937498488891
$ws.Range("R4").Value = 0.7550043748640001
$ws.Range("S4").Value = 0.005727817873891596
$ws.Range("T4").Value = 0.005727817873891595

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Alcam"
$ws.Range("C5").Value = "Cd6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5683613333333334
$ws.Range("H5").Value = 1.705084
$ws.Range("I5").Value = 0.008907897969731461
$ws.Range("J5").Value = 0.008907897969731461
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.08194666666666665
$ws.Range("N5").Value = 0.24584
$ws.Range("O5").Value = 0.356995568050465
$ws.Range("P5").Value = 0.356995568050465
$ws.Range("Q5").Value = 0.04657531672888889
$ws.Range("R5").Value = 0.41917785056
$ws.Range("S5").Value = 0.003180080095839866
$ws.Range("T5").Value = 0.003180080095839866

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Alcam"
$ws.Range("C6").Value = "Cd6"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.110262
$ws.Range("H6").Value = 3.330786
$ws.Range("I6").Value = 0.01740107927058724
$ws.Range("J6").Value = 0.01740107927058724
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1475986666666667
$ws.Range("N6").Value = 0.442796
$ws.Range("O6").Value = 0.6430044319495352
$ws.Range("P6").Value = 0.643004431949535
$ws.Range("Q6").Value = 0.1638731908506667
$ws.Range("R6").Value = 1.474858717656
$ws.Range("S6").Value = 0.01118897109169278
$ws.Range("T6").Value = 0.01118897109169278

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Alcam"
$ws.Range("C7").Value = "Cd6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.110262
$ws.Range("H7").Value = 3.330786
$ws.Range("I7").Value = 0.01740107927058724
$ws.Range("J7").Value = 0.01740107927058724
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.08194666666666665
$ws.Range("N7").Value = 0.24584
$ws.Range("O7").Value = 0.356995568050465
$ws.Range("P7").Value = 0.356995568050465
$ws.Range("Q7").Value = 0.09098227002666666
$ws.Range("R7").Value = 0.81884043024
$ws.Range("S7").Value = 0.006212108178894462
$ws.Range("T7").Value = 0.006212108178894464
